# Add two new enum-constrained columns ("integer_enum", "number_enum") to the
# "main" sheet, wire up their conditional formatting / data validation against
# two new lookup columns on the hidden "lists" sheet, and extend the existing
# "id" column's blank-row check to span the new columns.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("main")
$lists = $wb.Worksheets.Item("lists")

# ---------------------------------------------------------------------------
# 1. Populate the two new lookup columns on the hidden "lists" sheet.
# ---------------------------------------------------------------------------
$lists.Range("B1").Value = 1
$lists.Range("B2").Value = 2
$lists.Range("B3").Value = 3
$lists.Range("C1").Value = 1.25
$lists.Range("C2").Value = 1.5
$lists.Range("C3").Value = 1.75

# ---------------------------------------------------------------------------
# 2. New header cells on "main": H1 = integer_enum, I1 = number_enum.
#    Copy G1's formatting (bold header style) onto them first.
# ---------------------------------------------------------------------------
$main.Range("G1").Copy()
$main.Range("H1:I1").PasteSpecial(-4122) | Out-Null

$main.Range("H1").Value = "integer_enum"
$main.Range("I1").Value = "number_enum"

$hComment = $main.Range("H1").AddComment("One of the first three positive integers")
$iComment = $main.Range("I1").AddComment("One of the first three quarters after 1")

# Column widths for the new columns.
$main.Columns.Item(8).ColumnWidth = 14.3
$main.Columns.Item(9).ColumnWidth = 13

# ---------------------------------------------------------------------------
# 3. Conditional formatting.
#    - Extend column A's blank-row check to cover A:I (was A:G).
#    - Add new rules for H and I (enum membership against 'lists'!B / C).
# ---------------------------------------------------------------------------
$aCf = $main.Range("A2:A1048576").FormatConditions.Item(1)
$aCf.Formula1 = "=OR(AND(ISBLANK(A2), COUNTBLANK(`$A2:`$I2) <> 9), IF(ISBLANK(A2), FALSE, OR(IF(ISNUMBER(A2), INT(A2) <> A2, TRUE), COUNTIF(A`$2:A`$1048576, A2) >= 2, A2 < 1)))"

$hCf = $main.Range("H2:H1048576").FormatConditions.Add(2, 0, "=IF(ISBLANK(H2), FALSE, OR(IF(ISNUMBER(H2), INT(H2) <> H2, TRUE), ISNA(MATCH(H2, 'lists'!`$B`$1:`$B`$3, 0))))")
$hCf.Interior.Color = 13551615

$iCf = $main.Range("I2:I1048576").FormatConditions.Add(2, 0, "=IF(ISBLANK(I2), FALSE, OR(NOT(ISNUMBER(I2)), ISNA(MATCH(I2, 'lists'!`$C`$1:`$C`$3, 0))))")
$iCf.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# 4. Data validation: dropdown lists sourced from the new "lists" columns.
# ---------------------------------------------------------------------------
$hRange = $main.Range("H2:H1048576")
$hRange.Validation.Add(3, 3, 1, "='lists'!`$B`$1:`$B`$3") | Out-Null
$hDv = $hRange.Validation
$hDv.AlertStyle = 3
$hDv.ErrorTitle = "Invalid value"
$hDv.ErrorMessage = "Value must be in the dropdown list"
$hDv.IgnoreBlank = $true
$hDv.ShowInput = $true
$hDv.ShowError = $true

$iRange = $main.Range("I2:I1048576")
$iRange.Validation.Add(3, 3, 1, "='lists'!`$C`$1:`$C`$3") | Out-Null
$iDv = $iRange.Validation
$iDv.AlertStyle = 3
$iDv.ErrorTitle = "Invalid value"
$iDv.ErrorMessage = "Value must be in the dropdown list"
$iDv.IgnoreBlank = $true
$iDv.ShowInput = $true
$iDv.ShowError = $true

Write-Host "Done applying integer_enum / number_enum columns"
